# Weekly update: insert two new "Cebollín" price records (week of 2022-07-29,
# serial 44771) ahead of the existing history, shifting the rest of the
# table down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 28 (pushes old rows 28.. down to 30..)
$ws.Rows.Item(28).Insert()
$ws.Rows.Item(28).Insert()

# New row 28: "Primera" quality record for the new date
$ws.Cells.Item(28, 1).Value = 11
$ws.Cells.Item(28, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(28, 3).Value = "Bíobío"
$ws.Cells.Item(28, 4).Value = 44771
$ws.Cells.Item(28, 5).Value = 8
$ws.Cells.Item(28, 6).Value = 100112037
$ws.Cells.Item(28, 7).Value = "Cebollín"
$ws.Cells.Item(28, 8).Value = "Sin especificar"
$ws.Cells.Item(28, 9).Value = "Primera"
$ws.Cells.Item(28, 10).Value = 100
$ws.Cells.Item(28, 11).Value = 700
$ws.Cells.Item(28, 12).Value = 800
$ws.Cells.Item(28, 13).Value = 750
$ws.Cells.Item(28, 14).Value = "`$/paquete 6 unidades"
$ws.Cells.Item(28, 15).Value = "Región de Ñuble"
$ws.Cells.Item(28, 16).Value = 125
$ws.Cells.Item(28, 17).Value = 6
$ws.Cells.Item(28, 18).Value = "Hortaliza"

# New row 29: "Segunda" quality record for the same new date
$ws.Cells.Item(29, 1).Value = 11
$ws.Cells.Item(29, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(29, 3).Value = "Bíobío"
$ws.Cells.Item(29, 4).Value = 44771
$ws.Cells.Item(29, 5).Value = 8
$ws.Cells.Item(29, 6).Value = 100112037
$ws.Cells.Item(29, 7).Value = "Cebollín"
$ws.Cells.Item(29, 8).Value = "Sin especificar"
$ws.Cells.Item(29, 9).Value = "Segunda"
$ws.Cells.Item(29, 10).Value = 50
$ws.Cells.Item(29, 11).Value = 600
$ws.Cells.Item(29, 12).Value = 600
$ws.Cells.Item(29, 13).Value = 600
$ws.Cells.Item(29, 14).Value = "`$/paquete 6 unidades"
$ws.Cells.Item(29, 15).Value = "Región de Ñuble"
$ws.Cells.Item(29, 16).Value = 100
$ws.Cells.Item(29, 17).Value = 6
$ws.Cells.Item(29, 18).Value = "Hortaliza"
